$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the timestamp in column A for rows 2-11 to the new append time
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "2026-02-15 02:37:59"
}

# Swap the displayed URLs in F6 and F7 (underlying hyperlink targets stay as-is)
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5491569"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5491578"
